$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace the participant "Minna Hall" with "Frances Hamerstrom"
$ws.Range("A2").Value = "Frances Hamerstrom"

# Move the active selection to A2 (matches the new selection in the saved file)
$ws.Range("A2").Select()
